# Penalty Reward System (unfinished) - replicate the committed edit:
#  - Weekly Sales: append a new week (row 104)
#  - Daily PO: two cancelled POs (3T29EG9V / 3SFP4VCP, rows 43-44) are removed,
#    and the remaining PO rows 8-15/33-37/41-42 get re-keyed/re-valued to match
#    the new order
#  - Merged (Optional): the two corresponding merged-calendar rows (104-105,
#    both dated 2024-03-norm 45364) are removed and a new trailing week row
#    is appended; several Daily_PO_Qty (col C) values are corrected
#  - PO Volume Insights / PO Prediction: recomputed aggregate values

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Weekly Sales")
$ws2 = $wb.Worksheets.Item("Daily PO")
$ws3 = $wb.Worksheets.Item("Merged (Optional)")
$ws4 = $wb.Worksheets.Item("PO Volume Insights")
$ws5 = $wb.Worksheets.Item("PO Prediction")

# ---------------------------------------------------------------------------
# Weekly Sales: add new row 104 (new week, zero sales so far)
# ---------------------------------------------------------------------------
$ws1.Range("A104").NumberFormat = $ws1.Range("A103").NumberFormat
$ws1.Range("A104").Value = 45662.99999999999
$ws1.Range("B104").Value = 0

# ---------------------------------------------------------------------------
# Daily PO: update the re-keyed rows in place (values taken from the diff)
# ---------------------------------------------------------------------------
$ws2.Range("A8").Value = "1BZWCQBV"
$ws2.Range("M8").Value = 6
$ws2.Range("N8").Value = 6
$ws2.Range("P8").Value = 6
$ws2.Range("S8").Value = "LGB8"
$ws2.Range("T8").Value = 44965
$ws2.Range("U8").Value = 44972
$ws2.Range("W8").Value = 245
$ws2.Range("Y8").Value = 1470
$ws2.Range("Z8").Value = 1470
$ws2.Range("AA8").Value = 1470

$ws2.Range("A9").Value = "5I9ZRHZV"
$ws2.Range("S9").Value = "SBD1"

$ws2.Range("A10").Value = "4A7B7IQQ"
$ws2.Range("M10").Value = 12
$ws2.Range("N10").Value = 12
$ws2.Range("P10").Value = 12
$ws2.Range("S10").Value = "FTW1"
$ws2.Range("T10").Value = 45009
$ws2.Range("U10").Value = 45013
$ws2.Range("W10").Value = 244.99
$ws2.Range("Y10").Value = 2939.88
$ws2.Range("Z10").Value = 2939.88
$ws2.Range("AA10").Value = 2939.88

$ws2.Range("A11").Value = "7HZDRTMG"
$ws2.Range("M11").Value = 12
$ws2.Range("N11").Value = 12
$ws2.Range("P11").Value = 12
$ws2.Range("S11").Value = "FTW1"
$ws2.Range("T11").Value = 45041
$ws2.Range("U11").Value = 45043
$ws2.Range("Y11").Value = 2939.88
$ws2.Range("Z11").Value = 2939.88
$ws2.Range("AA11").Value = 2939.88

$ws2.Range("A13").Value = "7Y6ZLXLI"
$ws2.Range("M13").Value = 36
$ws2.Range("N13").Value = 36
$ws2.Range("P13").Value = 36
$ws2.Range("S13").Value = "LAX9"
$ws2.Range("T13").Value = 45009
$ws2.Range("U13").Value = 45013
$ws2.Range("Y13").Value = 8819.639999999999
$ws2.Range("Z13").Value = 8819.639999999999
$ws2.Range("AA13").Value = 8819.639999999999

$ws2.Range("A14").Value = "7MM8EEPT"
$ws2.Range("M14").Value = 24
$ws2.Range("N14").Value = 24
$ws2.Range("P14").Value = 24
$ws2.Range("S14").Value = "SCK4"
$ws2.Range("T14").Value = 45009
$ws2.Range("U14").Value = 45013
$ws2.Range("Y14").Value = 5879.76
$ws2.Range("Z14").Value = 5879.76
$ws2.Range("AA14").Value = 5879.76

$ws2.Range("A15").Value = "7UC1XMLE"
$ws2.Range("S15").Value = "SCK4"
$ws2.Range("T15").Value = 45041
$ws2.Range("U15").Value = 45043

$ws2.Range("A33").Value = "7GNF67EU"
$ws2.Range("M33").Value = 12
$ws2.Range("N33").Value = 12
$ws2.Range("P33").Value = 12
$ws2.Range("S33").Value = "ONT8"
$ws2.Range("Y33").Value = 2940
$ws2.Range("Z33").Value = 2940
$ws2.Range("AA33").Value = 2940

$ws2.Range("A34").Value = "7P58G3HL"
$ws2.Range("M34").Value = 6
$ws2.Range("N34").Value = 6
$ws2.Range("P34").Value = 6
$ws2.Range("S34").Value = "GYR3"
$ws2.Range("Y34").Value = 1470
$ws2.Range("Z34").Value = 1470
$ws2.Range("AA34").Value = 1470

$ws2.Range("A36").Value = "5ST9VQ1C"
$ws2.Range("M36").Value = 6
$ws2.Range("N36").Value = 6
$ws2.Range("P36").Value = 6
$ws2.Range("S36").Value = "GYR3"
$ws2.Range("Y36").Value = 1470
$ws2.Range("Z36").Value = 1470
$ws2.Range("AA36").Value = 1470

$ws2.Range("A37").Value = "4BFK6T1D"
$ws2.Range("M37").Value = 12
$ws2.Range("N37").Value = 12
$ws2.Range("P37").Value = 12
$ws2.Range("S37").Value = "ONT8"
$ws2.Range("Y37").Value = 2940
$ws2.Range("Z37").Value = 2940
$ws2.Range("AA37").Value = 2940

$ws2.Range("A41").Value = "79IMABND"
$ws2.Range("S41").Value = "SBD1"

$ws2.Range("A42").Value = "44NJRHWE"
$ws2.Range("S42").Value = "PHX6"

# Remove the two cancelled PO rows (previously rows 43 and 44: 3T29EG9V / 3SFP4VCP)
$ws2.Rows.Item(43).Delete()
$ws2.Rows.Item(43).Delete()

# ---------------------------------------------------------------------------
# Merged (Optional): fix Daily_PO_Qty (col C) values to match the corrected
# Daily PO sheet
# ---------------------------------------------------------------------------
$ws3.Range("C16").Value = 24
$ws3.Range("C17").Value = 36
$ws3.Range("C19").Value = 12
$ws3.Range("C21").Value = 6
$ws3.Range("C22").Value = 12
$ws3.Range("C66").Value = 12
$ws3.Range("C67").Value = 6
$ws3.Range("C70").Value = 6
$ws3.Range("C71").Value = 12

# Remove the two merged-calendar rows for the now-deleted POs (rows 104, 105)
$ws3.Rows.Item(104).Delete()
$ws3.Rows.Item(104).Delete()

# Append the new trailing week row (145) to mirror the new Weekly Sales week
$ws3.Range("A145").NumberFormat = $ws3.Range("A144").NumberFormat
$ws3.Range("A145").Value = 45662.99999999999
$ws3.Range("B145").Value = 0
$ws3.Range("C145").Value = 0

# ---------------------------------------------------------------------------
# PO Volume Insights: recomputed totals/average after removing the two
# cancelled POs
# ---------------------------------------------------------------------------
$ws4.Range("A2").Value = 390
$ws4.Range("B2").Value = 9.512195121951219

# ---------------------------------------------------------------------------
# PO Prediction: recomputed forecast value
# ---------------------------------------------------------------------------
$ws5.Range("A2").Value = 5.626829268292682
